$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Slitrk6"
$ws.Range("C2").Value = "Ptprs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.847512666666667
$ws.Range("H2").Value = 5.542538
$ws.Range("I2").Value = 0.8013421105175125
$ws.Range("J2").Value = 0.8013421105175124
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.127188333333333
$ws.Range("N2").Value = 12.381565
$ws.Range("O2").Value = 0.0561359176022362
$ws.Range("P2").Value = 0.05613591760223619
$ws.Range("Q2").Value = 7.625032723552222
$ws.Range("R2").Value = 68.62529451197
$ws.Range("S2").Value = 0.04498407468721313
$ws.Range("T2").Value = 0.04498407468721312

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Slitrk6"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.847512666666667
$ws.Range("H3").Value = 5.542538
$ws.Range("I3").Value = 0.8013421105175125
$ws.Range("J3").Value = 0.8013421105175124
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.24901333333333
$ws.Range("N3").Value = 141.74704
$ws.Range("O3").Value = 0.6426570597336346
$ws.Range("P3").Value = 0.6426570597336345
$ws.Range("Q3").Value = 87.29315062083556
$ws.Range("R3").Value = 785.6383555875201
$ws.Range("S3").Value = 0.5149881645859299
$ws.Range("T3").Value = 0.5149881645859297

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Slitrk6"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.847512666666667
$ws.Range("H4").Value = 5.542538
$ws.Range("I4").Value = 0.8013421105175125
$ws.Range("J4").Value = 0.8013421105175124
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.218847999999999
$ws.Range("N4").Value = 24.656544
$ws.Range("O4").Value = 0.1117885923419141
$ws.Range("P4").Value = 0.1117885923419141
$ws.Range("Q4").Value = 15.184425785408
$ws.Range("R4").Value = 136.659832068672
$ws.Range("S4").Value = 0.08958090651905126
$ws.Range("T4").Value = 0.08958090651905125

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slitrk6"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.847512666666667
$ws.Range("H5").Value = 5.542538
$ws.Range("I5").Value = 0.8013421105175125
$ws.Range("J5").Value = 0.8013421105175124
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.92629833333334
$ws.Range("N5").Value = 41.77889500000001
$ws.Range("O5").Value = 0.1894184303222152
$ws.Range("P5").Value = 0.1894184303222152
$ws.Range("Q5").Value = 25.72901257061223
$ws.Range("R5").Value = 231.5611131355101
$ws.Range("S5").Value = 0.1517889647253183
$ws.Range("T5").Value = 0.1517889647253183

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Slitrk6"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4580103333333334
$ws.Range("H6").Value = 1.374031
$ws.Range("I6").Value = 0.1986578894824876
$ws.Range("J6").Value = 0.1986578894824876
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.127188333333333
$ws.Range("N6").Value = 12.381565
$ws.Range("O6").Value = 0.0561359176022362
$ws.Range("P6").Value = 0.05613591760223619
$ws.Range("Q6").Value = 1.890294904279444
$ws.Range("R6").Value = 17.012654138515
$ws.Range("S6").Value = 0.01115184291502307
$ws.Range("T6").Value = 0.01115184291502307

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Slitrk6"
$ws.Range("C7").Value = "Ptprs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4580103333333334
$ws.Range("H7").Value = 1.374031
$ws.Range("I7").Value = 0.1986578894824876
$ws.Range("J7").Value = 0.1986578894824876
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 47.24901333333333
$ws.Range("N7").Value = 141.74704
$ws.Range("O7").Value = 0.6426570597336346
$ws.Range("P7").Value = 0.6426570597336345
$ws.Range("Q7").Value = 21.64053634647111
$ws.Range("R7").Value = 194.76482711824
$ws.Range("S7").Value = 0.1276688951477049
$ws.Range("T7").Value = 0.1276688951477048

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slitrk6"
$ws.Range("C8").Value = "Ptprs"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4580103333333334
$ws.Range("H8").Value = 1.374031
$ws.Range("I8").Value = 0.1986578894824876
$ws.Range("J8").Value = 0.1986578894824876
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.218847999999999
$ws.Range("N8").Value = 24.656544
$ws.Range("O8").Value = 0.1117885923419141
$ws.Range("P8").Value = 0.1117885923419141
$ws.Range("Q8").Value = 3.764317312096
$ws.Range("R8").Value = 33.878855808864
$ws.Range("S8").Value = 0.02220768582286283
$ws.Range("T8").Value = 0.02220768582286283

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slitrk6"
$ws.Range("C9").Value = "Ptprs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4580103333333334
$ws.Range("H9").Value = 1.374031
$ws.Range("I9").Value = 0.1986578894824876
$ws.Range("J9").Value = 0.1986578894824876
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 13.92629833333334
$ws.Range("N9").Value = 41.77889500000001
$ws.Range("O9").Value = 0.1894184303222152
$ws.Range("P9").Value = 0.1894184303222152
$ws.Range("Q9").Value = 6.378388541749445
$ws.Range("R9").Value = 57.40549687574501
$ws.Range("S9").Value = 0.03762946559689691
$ws.Range("T9").Value = 0.03762946559689691

